# Update "想去人数" (want-to-go count) figures for three events that
# appear in both the "展览" sheet and the combined "全部类型" sheet.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 189
$wsExhibit.Range("F5").Value = 106
$wsExhibit.Range("F6").Value = 626

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 189
$wsAll.Range("F7").Value = 106
$wsAll.Range("F8").Value = 626
